# Update data values as per algorithm re-run (KNN imputation) results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = -21.924
$ws.Range("C6").Value = -12.156
$ws.Range("D10").Value = -7.498
$ws.Range("A14").Value = -21.772
$ws.Range("D15").Value = -7.781999999999999
$ws.Range("C18").Value = -12.525
$ws.Range("C19").Value = -12.007
$ws.Range("A21").Value = -20.091
$ws.Range("D21").Value = -8.297999999999998
$ws.Range("B22").Value = 7.483
$ws.Range("D22").Value = -7.888000000000001
$ws.Range("A23").Value = -20.766
$ws.Range("B24").Value = 5.002
$ws.Range("D24").Value = -7.725999999999999
$ws.Range("A25").Value = -21.626
$ws.Range("A26").Value = -21.408
$ws.Range("B28").Value = 5.825
$ws.Range("A29").Value = -21.577
$ws.Range("D33").Value = -7.476000000000001
$ws.Range("B36").Value = 7.395
$ws.Range("C44").Value = -12.357
$ws.Range("B45").Value = 5.680000000000001
$ws.Range("D46").Value = -7.880000000000001
$ws.Range("C47").Value = -12.484
$ws.Range("B48").Value = 5.48
$ws.Range("B49").Value = 6.657999999999999
$ws.Range("D49").Value = -8.040000000000001
$ws.Range("C51").Value = -11.802
$ws.Range("B52").Value = 5.280000000000001
$ws.Range("A53").Value = -20.888
$ws.Range("B53").Value = 6.704000000000001
$ws.Range("B54").Value = 5.702
$ws.Range("C55").Value = -13.197
$ws.Range("D56").Value = -7.74
$ws.Range("A57").Value = -22.031
$ws.Range("C57").Value = -13.302
$ws.Range("A59").Value = -22.24
$ws.Range("D61").Value = -8.032
$ws.Range("C64").Value = -10.76
$ws.Range("D66").Value = -7.478999999999999
$ws.Range("A69").Value = -21.53
$ws.Range("B70").Value = 4.935
$ws.Range("D74").Value = -7.985000000000001
$ws.Range("D77").Value = -7.670999999999999
$ws.Range("A79").Value = -21.089
$ws.Range("C80").Value = -12.55
$ws.Range("A83").Value = -21.976
$ws.Range("B86").Value = 5.583
$ws.Range("B87").Value = 4.834999999999999
$ws.Range("D87").Value = -7.944
$ws.Range("D88").Value = -7.597000000000001
$ws.Range("B89").Value = 5.950000000000001
$ws.Range("A91").Value = -21.127
$ws.Range("C92").Value = -11.072
$ws.Range("A93").Value = -21.598
$ws.Range("C94").Value = -11.761
$ws.Range("C96").Value = -11.57
$ws.Range("D100").Value = -7.489999999999999
$ws.Range("B101").Value = 4.727
$ws.Range("C101").Value = -13.114
$ws.Range("A103").Value = -21.948
